$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.983.82'
$ws.Range("E2").Value = '  -1.15%  '
$ws.Range("D3").Value = '3.165.02'
$ws.Range("E3").Value = '  -4.49%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '591.41'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.61'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.30%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '3.160.51'
$ws.Range("E8").Value = '  -4.48%  '
$ws.Range("E9").Value = '  -0.99%  '
$ws.Range("E10").Value = '  -5.26%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.23'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.96%  '
$ws.Range("E12").Value = '  -3.00%  '
$ws.Range("E13").Value = '  -3.94%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.76'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.00%  '
$ws.Range("D15").Value = '3.688.38'
$ws.Range("E15").Value = '  -4.47%  '
$ws.Range("D17").Value = '3.167.53'
$ws.Range("E17").Value = '  -4.44%  '
$ws.Range("D18").Value = '62.915.21'
$ws.Range("E18").Value = '  -1.39%  '
$ws.Range("E19").Value = '  -4.39%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '460.45'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.06%  '
$ws.Range("E21").Value = '  -1.65%  '
$ws.Range("E22").Value = '  -4.93%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.61'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -6.46%  '
$ws.Range("E24").Value = '  -2.25%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.64'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.27%  '
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("E28").Value = '  -3.74%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.70'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.88%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.73'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -6.70%  '
$ws.Range("E31").Value = '  -5.81%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '27.10'
$ws.Range("D32").Style = "Normal"
$ws.Range("E33").Value = '  -2.60%  '
$ws.Range("E34").Value = '  -6.35%  '
$ws.Range("E35").Value = '  -6.69%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.81'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.02%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '51.14'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.00%  '
$ws.Range("D38").Value = '0.0₃0703'
$ws.Range("E38").Value = '  -5.00%  '
$ws.Range("E39").Value = '  -2.54%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '402.12'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -7.51%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.12'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.38%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.62'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.79%  '
$ws.Range("E43").Value = '  -7.15%  '
$ws.Range("D44").Value = '2.798.09'
$ws.Range("E44").Value = '  -9.70%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.251'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.00%  '
$ws.Range("E47").Value = '  -5.02%  '
$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '25.35'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.59%  '
$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '123.38'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.44%  '
$ws.Range("E50").Value = '  -2.06%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '34.25'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.77%  '
